$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'278.71"
$ws.Range("E2").Value = "'6.69%"
$ws.Range("D3").Value = "'27.17"
$ws.Range("E3").Value = "'0.55%"
$ws.Range("D4").Value = "'4.825"
$ws.Range("E4").Value = "'2.69%"
$ws.Range("D5").Value = "'0.06263"
$ws.Range("E5").Value = "'0.73%"
$ws.Range("D6").Value = "'6.861"
$ws.Range("E6").Value = "'1.62%"
$ws.Range("D7").Value = "'0.8786"
$ws.Range("E7").Value = "'3.17%"
$ws.Range("D8").Value = "'0.9434"
$ws.Range("E8").Value = "'2.94%"
$ws.Range("D9").Value = "'0.1450"
$ws.Range("E9").Value = "'3.32%"
$ws.Range("D10").Value = "'0.05145"
$ws.Range("E10").Value = "'6.57%"
$ws.Range("D11").Value = "'0.07282"
$ws.Range("E11").Value = "'2.82%"
$ws.Range("D12").Value = "'0.03161"
$ws.Range("E12").Value = "'1.73%"
$ws.Range("D13").Value = "'0.09051"
$ws.Range("E13").Value = "'-0.02%"
$ws.Range("D14").Value = "'0.001552"
$ws.Range("E14").Value = "'1.49%"
$ws.Range("D15").Value = "'0.0006276"
$ws.Range("E15").Value = "'1.64%"
$ws.Range("D16").Value = "'0.006057"
$ws.Range("E16").Value = "'0.20%"
$ws.Range("D17").Value = "'3.452"
$ws.Range("E17").Value = "'0.32%"
$ws.Range("D18").Value = "'3.263"
$ws.Range("E18").Value = "'2.79%"
$ws.Range("D19").Value = "'2.287"
$ws.Range("E19").Value = "'5.64%"
$ws.Range("E21").Value = "'-0.07%"
$ws.Range("D22").Value = "'3.853"
$ws.Range("E22").Value = "'-5.83%"
$ws.Range("D23").Value = "'0.04307"
$ws.Range("E23").Value = "'1.70%"
$ws.Range("E24").Value = "'-2.25%"
$ws.Range("D25").Value = "'0.004281"
$ws.Range("E25").Value = "'4.93%"
$ws.Range("E26").Value = "'-0.13%"
$ws.Range("D27").Value = "'0.0001688"
$ws.Range("E27").Value = "'2.94%"
$ws.Range("D40").Value = "'0.04040"
$ws.Range("E40").Value = "'2.21%"
$ws.Range("D41").Value = "'0.006410"
$ws.Range("E41").Value = "'55.67%"
$ws.Range("D42").Value = "'0.1154"
$ws.Range("E42").Value = "'3.74%"
$ws.Range("D43").Value = "'0.002104"
$ws.Range("E43").Value = "'-4.83%"
$ws.Range("D44").Value = "'0.01386"
$ws.Range("E44").Value = "'-0.16%"
$ws.Range("D45").Value = "'0.00005219"
$ws.Range("E45").Value = "'1.10%"
$ws.Range("E46").Value = "'-0.13%"
$ws.Range("D47").Value = "'2.352"
$ws.Range("E47").Value = "'1,007.83%"
$ws.Range("E49").Value = "'-0.13%"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("E50").Value = "'-0.13%"
